$wb = $excel.ActiveWorkbook

# 1. Status text changed from "Ready for handoff" to "In Translation"
#    (appears on the Overview sheet in the per-locale status columns, and
#    on each locale sheet's "Status" column).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # Force a string comparison - if $cell.Value2 is a [bool] (e.g. the
        # literal "True"/"False" cells elsewhere in the sheet), PowerShell's
        # "-eq" would coerce the right-hand string to bool instead and give
        # a false match, so stringify the left side first.
        if ("$($cell.Value2)" -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# 2. Narrow the per-locale status columns so the new (shorter) status text
#    doesn't need as much room.
#    Overview sheet: columns E (zh-cn) and F (de-de).
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

#    zh-cn / de-de sheets: column C (Status).
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
